$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows right after the header row (row 1), shifting existing
# data (old rows 2-21) down to rows 11-30.
$ws.Rows("2:10").Insert()

# The insert pulls in the header row's formatting; strip it so the new
# data rows stay unstyled like the rest of the numeric rows.
$ws.Range("A2:C10").ClearFormats()

# New data for the freshly inserted rows 2-10.
$newTop = @(
    @(0.0328340083360672, 0.0253509078174829, 0.0371100641787052),
    @(0.0198531206697225, 0.0163406450301408, -0.009468411095440299),
    @(-0.0001527163112768, 0.0464257597923278, 0.0154243474826216),
    @(0.0108428578823804, -0.015118914656341, 0.1504255682229995),
    @(0.0204639863222837, -0.0316122770309448, 0.0610865242779254),
    @(-0.009010262787342, -0.0128281703218817, 0.0167987942695617),
    @(-0.0468839071691036, 0.052381694316864, 0.0694859251379966),
    @(0.0027488935738801, 0.09178250283002851, 0.0717766657471656),
    @(0.0128281703218817, 0.0387899428606033, 0.0054977871477603)
)

for ($i = 0; $i -lt $newTop.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $newTop[$i][0]
    $ws.Cells.Item($r, 2).Value = $newTop[$i][1]
    $ws.Cells.Item($r, 3).Value = $newTop[$i][2]
}

# Append a brand new row of data at the end (row 31).
$ws.Cells.Item(31, 1).Value = 0.3949243724346161
$ws.Cells.Item(31, 2).Value = 0.0597120784223079
$ws.Cells.Item(31, 3).Value = 0.08445212244987479
